# "beetle 2 boss added"
#
# The "BOSS" section of Sheet1 had leftover placeholder text (Beatles
# member names / "Pete Best") in the description column for the
# "Prideful Worm" and "boss caterpillar" rows. This change:
#   - fixes the "Prideful Worm" description to end with a period, and
#   - replaces the stray "Pete Best" placeholder with the real
#     description for the new "beetle 2" boss row.
# It also updates the saved window scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B44: "Prideful Worm" row description -> add trailing period.
$ws.Range("B44").Value = "His Pride exceeds the limits for a worm."

# B45: "boss caterpillar" row description -> replace leftover
# "Pete Best" placeholder with the real "beetle 2" description.
$ws.Range("B45").Value = "I'm beautiful enough even without metamorphosis"

# Update the view: scroll position and active selection.
$ws.Activate() | Out-Null
$ws.Range("B46").Select() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 25
$aw.ScrollColumn = 1
